$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Auto-generated cell updates
# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "20.377.06"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.14%  "
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.460.42"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.28%  "
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.010"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.94%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9498"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.07%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "274.66"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.48%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3646"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.58%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3068"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.04%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "39.65"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.31%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.033"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.07%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06571"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.65%  "
# Row 12
$ws.Range("E12").Value = "  +0.05%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.418"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.19%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.83"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.29%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.125"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.06%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001023"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.53%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.460.32"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.27%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9679"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.23%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.05791"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.27%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.53"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.30%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.432"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.28%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.41"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.17%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.87"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.34%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.249"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.56%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "20.415.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.28%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "141.48"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +6.43%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.077"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -8.19%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.10"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.92%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.615.98"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.89%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "111.97"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.12%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.820"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.93%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.884"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.49%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.07877"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.39%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7887"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.84%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.525"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.71%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05702"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.46%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.144"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.46%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.669"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.01%  "
# Row 39
$ws.Range("B39").Value = "Frax"
$ws.Range("C39").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9573"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.00%  "
# Row 40
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02022"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.86%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "10.30"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.46%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.469"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -10.39%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1854"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.59%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5251"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.09%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.486"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.45%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "11.84"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.38%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "116.96"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.09%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5113"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.40%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.747"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.15%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06411"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.60%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9881"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.21%  "
